# A-I_Classifier_Modes_Transport.xlsx edit
#
# Intent (per commit message + xml diff): the "Plain English" column header
# used across the *_to_Code lookup sheets is renamed to "Plain_English"
# (underscore instead of space) so it is a valid key/identifier when the
# workbook's codes get consumed elsewhere (e.g. written out to YAML).
#
# Because this string is shared (sharedStrings.xml) across four sheets,
# retyping the header cell on any one of them updates it everywhere the
# same shared string is referenced - matching the diff, which shows every
# "Plain English" header cell (Fuel_to_Code!C1, VehFuel_to_Code!C1,
# Tech_to_Code!C1, Dem_to_Code!B1) switching from the old shared string to
# the new "Plain_English" one.
#
# The rest of the session also leaves the workbook parked on the
# Dem_to_Code sheet (the diff shows tabSelected/activeTab moving there,
# away from Fuel_per_VehFuel), with the specific last-used selections shown
# in the diff for the sheets that were visited along the way.

$wb = $excel.ActiveWorkbook

# --- Rename the "Plain English" header to "Plain_English" everywhere ---
$ws = $wb.Worksheets.Item("Fuel_to_Code")
[void]($ws.Range("C1").Value = "Plain_English")

$ws = $wb.Worksheets.Item("VehFuel_to_Code")
[void]($ws.Range("C1").Value = "Plain_English")

$ws = $wb.Worksheets.Item("Tech_to_Code")
[void]($ws.Range("C1").Value = "Plain_English")

$ws = $wb.Worksheets.Item("Dem_to_Code")
[void]($ws.Range("B1").Value = "Plain_English")

# --- Recreate the final navigation / selection state from the diff ---
$ws = $wb.Worksheets.Item("Fuel_per_VehFuel")
[void]$ws.Activate()
[void]$ws.Range("A10:XFD11").Select()

$ws = $wb.Worksheets.Item("Fuel_to_Code")
[void]$ws.Activate()
[void]$ws.Range("D8").Select()

$ws = $wb.Worksheets.Item("VehFuel_to_Code")
[void]$ws.Activate()
[void]$ws.Range("A10:XFD11").Select()

$ws = $wb.Worksheets.Item("Tech_to_Code")
[void]$ws.Activate()
[void]$ws.Range("A1:D13").Select()

$ws = $wb.Worksheets.Item("Dem_to_Code")
[void]$ws.Activate()
[void]$ws.Range("F13").Select()
